# Steffens tidsregistrering fra den 8/3
# Adds Steffen's time entries for 7/3 (date already has a row, just needed
# the date stamp) and fills in the 8/3 block (rows 11-14) with activities,
# start/end times and role counts. Also updates the sheet's selection to
# reflect where the user ended up (E15) after entering the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tuesday 7/3 block: stamp the date (activities/times already present) ---
$ws.Range("A7").Value = 42801

# --- Wednesday 8/3 block: stamp the date ---
$ws.Range("A10").Value = 42802

# New activity rows for 8/3
$ws.Range("F11").Value = "krydstjek"
$ws.Range("G11").Value = 0.34375
$ws.Range("H11").Value = 0.38541666666666669
$ws.Range("I11").Value = 1

$ws.Range("F12").Value = "review"
$ws.Range("G12").Value = 0.39583333333333331
$ws.Range("H12").Value = 0.47916666666666669
$ws.Range("I12").Value = 2

$ws.Range("F13").Value = "spørgsmål"
$ws.Range("G13").Value = 0.5
$ws.Range("H13").Value = 0.54166666666666663
$ws.Range("I13").Value = 1

$ws.Range("F14").Value = "OC 1-3-4"
$ws.Range("G14").Value = 0.54166666666666663
$ws.Range("H14").Value = 0.64583333333333337
$ws.Range("I14").Value = 3

# Start/end-time columns use the workbook's existing time format (h:mm)
$ws.Range("G11:H14").NumberFormat = "h:mm"

# Leave the selection where the author ended up after entering the data
$ws.Range("E15").Select()
